# Update countries & provincias Spain
# Applies the "after" edit described by the OOXML diff:
#  - Swap the country labels for 4 pairs of rows (the underlying data for each
#    country moved to a new row while being refreshed with newer numbers)
#  - Refresh the numeric statistics for a handful of rows
#  - Update the "Datos actualizados" timestamp in cell A1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp text (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 16:22"

# --- Row 15 (Brasil) : refreshed counts ---
$ws.Range("B15").Value = 34485
$ws.Range("C15").Value = 803
$ws.Range("E15").Value = 18278
$ws.Range("G15").Value = 40
$ws.Range("H15").Value = 2181

# --- Row 65 (Islandia) : refreshed counts ---
$ws.Range("B65").Value = 1760
$ws.Range("C65").Value = 6
$ws.Range("D65").Value = 1291
$ws.Range("E65").Value = 460
$ws.Range("F65").Value = 3

# --- Rows 107/108 : Jordania & Reunion swap places, Reunion gets new numbers ---
$ws.Range("A107").Value = "Reunion"
$ws.Range("C107").Value = 5
$ws.Range("D107").Value = 237
$ws.Range("E107").Value = 170
$ws.Range("F107").Value = 4
$ws.Range("H107").Value = 0

$ws.Range("A108").Value = "Jordania"
$ws.Range("B108").Value = 407
$ws.Range("D108").Value = 265
$ws.Range("E108").Value = 135
$ws.Range("F108").Value = 5
$ws.Range("H108").Value = 7

# --- Row 112 (Mauricio) : refreshed counts ---
$ws.Range("B112").Value = 325
$ws.Range("C112").Value = 1
$ws.Range("E112").Value = 208

# --- Rows 118/119 : Sri Lanka & Mayotte swap places, Mayotte gets new numbers ---
$ws.Range("A118").Value = "Mayotte"
$ws.Range("B118").Value = 254
$ws.Range("C118").Value = 9
$ws.Range("D118").Value = 117
$ws.Range("E118").Value = 133
$ws.Range("F118").Value = 6
$ws.Range("H118").Value = 4

$ws.Range("A119").Value = "Sri Lanka"
$ws.Range("B119").Value = 248
$ws.Range("C119").Value = 4
$ws.Range("D119").Value = 86
$ws.Range("E119").Value = 155
$ws.Range("F119").Value = 1
$ws.Range("H119").Value = 7

# --- Rows 122/123 : Paraguay & Mali swap places, Mali gets new numbers ---
$ws.Range("A122").Value = "Mali"
$ws.Range("B122").Value = 216
$ws.Range("C122").Value = 45
$ws.Range("D122").Value = 41
$ws.Range("E122").Value = 162
$ws.Range("F122").Value = 0
$ws.Range("H122").Value = 13

$ws.Range("A123").Value = "Paraguay"
$ws.Range("B123").Value = 202
$ws.Range("C123").Value = 3
$ws.Range("D123").Value = 35
$ws.Range("E123").Value = 159
$ws.Range("F123").Value = 1
$ws.Range("H123").Value = 8

# --- Rows 173/174 : Nepal & Sierra Leona swap places, Sierra Leona gets new numbers ---
$ws.Range("A173").Value = "Sierra Leona"
$ws.Range("C173").Value = 4
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 30

$ws.Range("A174").Value = "Nepal"
$ws.Range("B174").Value = 30
$ws.Range("D174").Value = 2
$ws.Range("E174").Value = 28
